$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'71.062.29"
$ws.Range("E2").Value2 = "  +6.13%  "

$ws.Range("D3").Value2 = "'3.661.64"
$ws.Range("E3").Value2 = "  +6.11%  "

$ws.Range("E4").Value2 = "  +0.03%  "

$ws.Range("D5").Value2 = "'596.87"
$ws.Range("E5").Value2 = "  +2.94%  "

$ws.Range("D6").Value2 = "'194.89"
$ws.Range("E6").Value2 = "  +3.20%  "

$ws.Range("D7").Value2 = "'0.648"
$ws.Range("E7").Value2 = "  +2.56%  "

$ws.Range("D8").Value2 = "'3.654.22"
$ws.Range("E8").Value2 = "  +6.13%  "

$ws.Range("E9").Value2 = "  +0.00%  "

$ws.Range("E10").Value2 = "  +7.83%  "

$ws.Range("D11").Value2 = "'0.676"
$ws.Range("E11").Value2 = "  +4.80%  "

$ws.Range("D12").Value2 = "'58.28"
$ws.Range("E12").Value2 = "  +2.39%  "

$ws.Range("D13").Value2 = "'0.0000296"
$ws.Range("E13").Value2 = "  +6.39%  "

$ws.Range("D14").Value2 = "'9.99"
$ws.Range("E14").Value2 = "  +5.71%  "

$ws.Range("D15").Value2 = "'4.243.78"

$ws.Range("D16").Value2 = "'20.18"
$ws.Range("E16").Value2 = "  +7.34%  "

$ws.Range("D17").Value2 = "'3.657.89"
$ws.Range("E17").Value2 = "  +6.33%  "

$ws.Range("D18").Value2 = "'71.036.23"
$ws.Range("E18").Value2 = "  +6.26%  "

$ws.Range("E19").Value2 = "  +6.01%  "

$ws.Range("E20").Value2 = "  +1.99%  "

$ws.Range("E21").Value2 = "  +4.18%  "

$ws.Range("D22").Value2 = "'489.62"
$ws.Range("E22").Value2 = "  +1.39%  "

$ws.Range("D23").Value2 = "'19.13"
$ws.Range("E23").Value2 = "  +12.79%  "

$ws.Range("D24").Value2 = "'5.28"
$ws.Range("E24").Value2 = "  -1.12%  "

$ws.Range("E25").Value2 = "  +3.94%  "

$ws.Range("D26").Value2 = "'91.44"
$ws.Range("E26").Value2 = "  +2.25%  "

$ws.Range("D27").Value2 = "'3.17"
$ws.Range("E27").Value2 = "  +6.41%  "

$ws.Range("D28").Value2 = "'11.49"
$ws.Range("E28").Value2 = "  +4.74%  "

$ws.Range("D29").Value2 = "'9.64"
$ws.Range("E29").Value2 = "  +6.36%  "

$ws.Range("D30").Value2 = "'32.93"
$ws.Range("E30").Value2 = "  +5.31%  "

$ws.Range("E31").Value2 = "  +5.64%  "

$ws.Range("E32").Value2 = "  +9.49%  "

$ws.Range("D33").Value2 = "'630.02"
$ws.Range("E33").Value2 = "  +5.16%  "

$ws.Range("D34").Value2 = "'12.31"
$ws.Range("E34").Value2 = "  +4.42%  "

$ws.Range("D35").Value2 = "'66.49"
$ws.Range("E35").Value2 = "  +3.64%  "

$ws.Range("D36").Value2 = "'40.33"
$ws.Range("E36").Value2 = "  +9.42%  "

$ws.Range("D37").Value2 = "'0.0₃0834"
$ws.Range("E37").Value2 = "  +10.55%  "

$ws.Range("E38").Value2 = "  +6.23%  "

$ws.Range("E39").Value2 = "  +0.10%  "

$ws.Range("E40").Value2 = "  +0.00%  "

$ws.Range("E41").Value2 = "  +2.16%  "

$ws.Range("D42").Value2 = "'3.335.56"

$ws.Range("D43").Value2 = "'3.16"
$ws.Range("E43").Value2 = "  +8.90%  "

$ws.Range("D44").Value2 = "'2.83"
$ws.Range("E44").Value2 = "  +11.68%  "

$ws.Range("D45").Value2 = "'3.10"
$ws.Range("E45").Value2 = "  +9.34%  "

$ws.Range("E46").Value2 = "  +5.70%  "

$ws.Range("D47").Value2 = "'9.44"
$ws.Range("E47").Value2 = "  +9.12%  "

$ws.Range("B48").Value2 = "ApeXProtocol"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value2 = "'3.33"
$ws.Range("E48").Value2 = "  +2.57%  "

$ws.Range("B49").Value2 = "Stellar"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value2 = "'0.140"
$ws.Range("E49").Value2 = "  +3.74%  "

$ws.Range("D50").Value2 = "'3.25"
$ws.Range("E50").Value2 = "  -3.21%  "

$ws.Range("E51").Value2 = "  +0.12%  "
